$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = "8 presses"
$ws.Range("D5").Value = "8 presses"

$ws.Columns.Item(4).ColumnWidth = 20.65

$ws.Range("D6").Select()
